# The sheet gained one new data row: a new record was inserted right
# before the existing row 700 (a "Primera" quality Pepino ensalada entry
# for Región de Arica y Parinacota, $/caja 60 unidades), pushing every
# row from 700 through 802 down by one (to 701 through 803).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 700; everything below shifts down.
$ws.Rows.Item(700).Insert()

# Populate the newly inserted row 700 with the new record's data.
$ws.Range("A700").Value = 8
$ws.Range("B700").Value = "Terminal La Palmera de La Serena"
$ws.Range("C700").Value = "Coquimbo"
$ws.Range("D700").Value = 44984
$ws.Range("E700").Value = 4
$ws.Range("F700").Value = 100112043
$ws.Range("G700").Value = "Pepino ensalada"
$ws.Range("H700").Value = "Sin especificar"
$ws.Range("I700").Value = "Primera"
$ws.Range("J700").Value = 400
$ws.Range("K700").Value = 8000
$ws.Range("L700").Value = 9000
$ws.Range("M700").Value = 8500
$ws.Range("N700").Value = '$/caja 60 unidades'
$ws.Range("O700").Value = "Región de Arica y Parinacota"
$ws.Range("P700").Value = 142
$ws.Range("Q700").Value = 60
$ws.Range("R700").Value = "Hortaliza"
